$d = $word.ActiveDocument

# 1. Title: "QRy" -> "QuinceQR" (rest of the heading stays identical)
$d.Content.Find.Execute("QRy", $false, $false, $false, $false, $false, $true, 1, $false, "QuinceQR", 2) | Out-Null

# 2. Relocate the "_GoBack" bookmark: it now sits at the very end of the
#    "Git & GitHub" bullet (it previously sat at the end of the "Bilder..."
#    bullet). Re-adding a bookmark with the same name moves it, since Word
#    only ever keeps a single "_GoBack" bookmark.
$find = $d.Content.Find
$found = $find.Execute("GitHub")
if ($found) {
    $gitHubRange = $find.Parent
    $d.Bookmarks.Add("_GoBack", $gitHubRange) | Out-Null
}
